# Update ticket/interest counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 138
$ws1.Range("F3").Value = 453
$ws1.Range("F4").Value = 4
$ws1.Range("F5").Value = 15
$ws1.Range("F7").Value = 26
$ws1.Range("F8").Value = 6
$ws1.Range("F9").Value = 74
$ws1.Range("G9").Value = "不可售"

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 138
$ws4.Range("F3").Value = 76
$ws4.Range("F4").Value = 453
$ws4.Range("F5").Value = 4
$ws4.Range("F6").Value = 15
$ws4.Range("F8").Value = 26
$ws4.Range("F9").Value = 6
$ws4.Range("F10").Value = 74
$ws4.Range("G10").Value = "不可售"
